# Case_6_51 / lines_states.xlsx update ("contingencies with rene fine")
# Inserts two new line rows (line7, line8) into the lines/contingencies table,
# pushing the existing extrN rows down by two, and refreshes the C/D/E values
# for the shifted + new rows to match the updated contingency run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert two blank rows right before the old row 8 (extr1) ---
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# --- 2. Give column A of the two new rows the same look (bold/border/centered)
#        as the rest of the index column ---
foreach ($r in 8, 9) {
    $idxCell = $ws.Cells.Item($r, 1)
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
}

# --- 3. Full target state for rows 8-17 (A:index, B:name, C:from_bus, D:to_bus, E:in_service) ---
$rows = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $false },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
